$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) hold numeric-looking text (e.g. "26.220.33", "0.0613")
# in the source data; force text format so COM does not coerce them to
# floating point numbers (which would round-trip with FP noise).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.220.33'
$ws.Range('E2').Value = '  +0.70%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.605.26'
$ws.Range('E3').Value = '  +0.51%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.49'
$ws.Range('E5').Value = '  -0.01%  '
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.485'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +0.31%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0613'
$ws.Range('E9').Value = '  -0.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.19'
$ws.Range('E10').Value = '  +1.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0814'
$ws.Range('E11').Value = '  -0.33%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.828.30'
$ws.Range('E12').Value = '  +0.49%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.626.13'
$ws.Range('E13').Value = '  +1.72%  '
$ws.Range('E14').Value = '  +0.36%  '
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '26.192.42'
$ws.Range('E17').Value = '  +2.15%  '
$ws.Range('E18').Value = '  +0.95%  '
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '200.55'
$ws.Range('E20').Value = '  -1.46%  '
$ws.Range('E21').Value = '  +0.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.27'
$ws.Range('E22').Value = '  -0.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.86'
$ws.Range('E24').Value = '  +1.72%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.04'
$ws.Range('E25').Value = '  +1.93%  '
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('E27').Value = '  -2.95%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.18'
$ws.Range('E28').Value = '  -0.25%  '
$ws.Range('E29').Value = '  +1.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0488'
$ws.Range('E30').Value = '  +3.73%  '
$ws.Range('E31').Value = '  +0.44%  '
$ws.Range('E32').Value = '  +2.38%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.93'
$ws.Range('E33').Value = '  -0.74%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.48'
$ws.Range('E34').Value = '  +0.68%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.39'
$ws.Range('E35').Value = '  +1.49%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.164.49'
$ws.Range('E36').Value = '  +4.95%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0167'
$ws.Range('E37').Value = '  +3.77%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.31'
$ws.Range('E39').Value = '  -0.51%  '
$ws.Range('E40').Value = '  +1.01%  '
$ws.Range('E41').Value = '  +0.75%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.784'
$ws.Range('E42').Value = '  +1.09%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.29'
$ws.Range('E43').Value = '  +3.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.740.49'
$ws.Range('E44').Value = '  +0.48%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '91.69'
$ws.Range('E45').Value = '  -0.74%  '
$ws.Range('E46').Value = '  +2.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '54.03'
$ws.Range('E47').Value = '  +1.27%  '
$ws.Range('E48').Value = '  +0.11%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0₇0972'
$ws.Range('E49').Value = '  +3.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.407'
$ws.Range('E50').Value = '  -0.51%  '
$ws.Range('E51').Value = '  -0.07%  '
